$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the crypto listing rows.
# D-column values are forced to text format ("@") so that Excel does not
# auto-convert numeric-looking strings (e.g. "20.00", "0.0000188") into
# actual numbers and strip formatting/precision.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "75.713.61"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.889.81"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "195.92"
$ws.Range("E5").Value = "  +2.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "593.86"
$ws.Range("E6").Value = "  -2.35%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.191"
$ws.Range("E9").Value = "  -4.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.887.89"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("E11").Value = "  +9.94%  "

$ws.Range("E12").Value = "  -1.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.86"
$ws.Range("E13").Value = "  -1.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.434.83"
$ws.Range("E14").Value = "  +1.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.608.28"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.18"
$ws.Range("E16").Value = "  -2.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000188"
$ws.Range("E17").Value = "  -3.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.895.25"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.82"
$ws.Range("E19").Value = "  -5.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.56"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.60"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.27"
$ws.Range("E22").Value = "  -2.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("E23").Value = "  -1.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.15"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.024.29"
$ws.Range("E26").Value = "  +2.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.17"
$ws.Range("E27").Value = "  -3.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  -3.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  -1.13%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.39"
$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "498.33"
$ws.Range("E32").Value = "  -8.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.66"
$ws.Range("E33").Value = "  -4.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -3.16%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.26"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.00"
$ws.Range("E37").Value = "  -3.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.70"
$ws.Range("E38").Value = "  +2.02%  "

$ws.Range("E39").Value = "  -8.32%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "179.37"
$ws.Range("E41").Value = "  -5.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.342"
$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.96"
$ws.Range("E43").Value = "  -5.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.65"
$ws.Range("E44").Value = "  -5.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0899"
$ws.Range("E45").Value = "  +4.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").Value = "  -6.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.06"
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -5.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.572"
$ws.Range("E49").Value = "  -2.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.70"
$ws.Range("E50").Value = "  -3.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.648"
$ws.Range("E51").Value = "  +4.87%  "

